$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
# New column layout (A..N):
# A Name | B Total Amount | C Value | D RAI | E Admit Score | F Major |
# G ACT Math | H ACT English | I ACT Composite | J SAT Math | K SAT Reading |
# L SAT Combined | M GPA | N HS Percentile

$ws.Range("D1").Value = "RAI"
$ws.Range("E1").Value = "Admit Score"
$ws.Range("F1").Value = "Major"
$ws.Range("G1").Value = "ACT Math"
$ws.Range("H1").Value = "ACT English"
$ws.Range("I1").Value = "ACT Composite"
$ws.Range("J1").Value = "SAT Math"
$ws.Range("K1").Value = "SAT Reading"
$ws.Range("L1").Value = "SAT Combined"
$ws.Range("M1").Value = "GPA"
$ws.Range("N1").Value = "HS Percentile"

# Carry the header style (bold/centered/bordered, thin box border) to the
# newly added header cells so they match the existing header formatting.
# Use copy/paste-special-formats (rather than setting Font/Alignment/Borders
# individually) so the new cells reuse the existing style record instead of
# generating near-duplicate style entries.
$ws.Range("C1").Copy()
$ws.Range("D1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2 (shift existing values into the new layout, add new data) ---
$ws.Range("D2").Value = 315
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = "All"
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = 27
$ws.Range("I2").Value = 26
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 96

# --- Row 3 (new "Test Two" record) --------------------------------------
$ws.Range("A3").Value = "Test Two"

# B3/C3 are text (numeric-looking strings), not numbers. Writing "50" /
# "10000" straight to Value would be auto-coerced to a number, so stage the
# text in a scratch cell formatted as Text, then copy/paste-special *values
# only* into B3/C3 - that brings over the text type without attaching any
# number-format style to the destination cell (matches the un-styled target
# cells). The scratch cells are cleared afterwards so they don't affect the
# sheet's used range/dimension.
$scratch1 = $ws.Range("Z99")
$scratch1.NumberFormat = "@"
$scratch1.Value = "50"
$scratch1.Copy()
$ws.Range("B3").PasteSpecial(-4163)

$scratch2 = $ws.Range("Z100")
$scratch2.NumberFormat = "@"
$scratch2.Value = "10000"
$scratch2.Copy()
$ws.Range("C3").PasteSpecial(-4163)

$ws.Range("Z99:Z100").Clear()
$excel.CutCopyMode = $false

$ws.Range("D3").Value = 330
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = "All"
$ws.Range("G3").Value = 32
$ws.Range("H3").Value = 28
$ws.Range("I3").Value = 30
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 620
$ws.Range("L3").Value = 1320
$ws.Range("M3").Value = 3.8
$ws.Range("N3").Value = 90
